# Apply the scheduled Kraken_Profits data refresh across all Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 751.5
$ws.Range("J12").Value = 751.5
$ws.Range("L12").Value = 751.5
$ws.Range("N12").Value = -1091.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 19414.143
$ws.Range("J70").Value = 21816.5
$ws.Range("L70").Value = 65449.5
$ws.Range("N70").Value = -65989.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 19414.143
$ws.Range("J73").Value = 21816.5
$ws.Range("L73").Value = 65449.5
$ws.Range("N73").Value = -67321.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1749.8889
$ws.Range("I98").Value = 1328.4286
$ws.Range("J98").Value = 3225
$ws.Range("K98").Value = 1328.4286
$ws.Range("L98").Value = 3225
$ws.Range("M98").Value = 169.5714
$ws.Range("N98").Value = -6221

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1749.8889
$ws.Range("I122").Value = 1328.4286
$ws.Range("J122").Value = 3225
$ws.Range("K122").Value = 3985.2858
$ws.Range("L122").Value = 9675
$ws.Range("M122").Value = -1535.2858
$ws.Range("N122").Value = -14575

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 250398.75
$ws.Range("I131").Value = 250398.75
$ws.Range("K131").Value = 751196.25
$ws.Range("M131").Value = -746156.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2220.625
$ws.Range("I137").Value = 2252.1428
$ws.Range("K137").Value = 6756.428400000001
$ws.Range("M137").Value = -4206.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3915.5
$ws.Range("I26").Value = 3624.75
$ws.Range("J26").Value = 4497
$ws.Range("K26").Value = 3624.75
$ws.Range("L26").Value = 4497
$ws.Range("M26").Value = -3294.75
$ws.Range("N26").Value = -5157

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3201
$ws.Range("I32").Value = 3334.4443
$ws.Range("K32").Value = 3334.4443
$ws.Range("M32").Value = -3047.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 499.57144
$ws.Range("I80").Value = 323.75
$ws.Range("J80").Value = 734
$ws.Range("K80").Value = 323.75
$ws.Range("L80").Value = 734
$ws.Range("M80").Value = 674.25
$ws.Range("N80").Value = -2730

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 499.57144
$ws.Range("I83").Value = 323.75
$ws.Range("J83").Value = 734
$ws.Range("K83").Value = 1618.75
$ws.Range("L83").Value = 3670
$ws.Range("M83").Value = 3373.25
$ws.Range("N83").Value = -13654

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 44800
$ws.Range("J48").Value = 44800
$ws.Range("L48").Value = 44800
$ws.Range("N48").Value = -45752

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 99995
$ws.Range("J100").Value = 99995
$ws.Range("L100").Value = 99995
$ws.Range("N100").Value = -102159

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3250
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 646.4
$ws.Range("I7").Value = 646.4
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1939.2
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1827.2
$ws.Range("N7").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2598.2144
$ws.Range("I103").Value = 2198.8
$ws.Range("J103").Value = 2820.111
$ws.Range("K103").Value = 6596.400000000001
$ws.Range("L103").Value = 8460.332999999999
$ws.Range("M103").Value = -5717.400000000001
$ws.Range("N103").Value = -10218.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 564.75
$ws.Range("I114").Value = 364
$ws.Range("K114").Value = 1092
$ws.Range("M114").Value = 2162

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 465.69232
$ws.Range("I117").Value = 375.8
$ws.Range("J117").Value = 521.875
$ws.Range("K117").Value = 1127.4
$ws.Range("L117").Value = 1565.625
$ws.Range("M117").Value = 2314.6
$ws.Range("N117").Value = -8449.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2359.625
$ws.Range("I129").Value = 1500
$ws.Range("J129").Value = 2646.1667
$ws.Range("K129").Value = 4500
$ws.Range("L129").Value = 7938.500100000001
$ws.Range("M129").Value = 500
$ws.Range("N129").Value = -17938.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1666.5714
$ws.Range("I97").Value = 1744.3334
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 1744.3334
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -1248.3334
$ws.Range("N97").Value = -2192

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 25026.5
$ws.Range("J42").Value = 20028
$ws.Range("L42").Value = 20028
$ws.Range("N42").Value = -21154

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 21537.334
$ws.Range("I43").Value = 7012
$ws.Range("J43").Value = 28800
$ws.Range("K43").Value = 7012
$ws.Range("L43").Value = 28800
$ws.Range("M43").Value = -6819
$ws.Range("N43").Value = -29186

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 25026.5
$ws.Range("J49").Value = 20028
$ws.Range("L49").Value = 20028
$ws.Range("N49").Value = -20322

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 16000
$ws.Range("J50").Value = 12000
$ws.Range("L50").Value = 12000
$ws.Range("N50").Value = -13274

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -40990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3097
$ws.Range("I136").Value = 3097
$ws.Range("K136").Value = 9291
$ws.Range("M136").Value = -6741

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 53000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 40301
$ws.Range("J80").Value = 40301
$ws.Range("L80").Value = 40301
$ws.Range("N80").Value = -42297

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 40301
$ws.Range("J83").Value = 40301
$ws.Range("L83").Value = 120903
$ws.Range("N83").Value = -130887

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = ""
$ws.Range("N96").Value = ""
